$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(1, 8).Value2 = "SSE"
$ws.Cells.Item(1, 9).Value2 = "SZSE"

$splitValues = @{
    2 = 411950
    3 = 1154954
    4 = 22557
    5 = 6128
}

$numFmt = '_(* #,##0_);_(* \(#,##0\);_(* "-"??_);_(@_)'

foreach ($r in 2, 3, 4, 5) {
    $old = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 8).Value2 = $old
    $ws.Cells.Item($r, 9).Value2 = $splitValues[$r]
    $ws.Cells.Item($r, 8).NumberFormat = $numFmt
    $ws.Cells.Item($r, 9).NumberFormat = $numFmt
    $ws.Cells.Item($r, 3).Formula = "=SUM(H" + $r + ":I" + $r + ")"
}

# Reposition / resize the chart to its new anchor
$co = $ws.ChartObjects(1)
$co.Left = 642.58984375
$co.Top = 110.99992125984252
$co.Width = 561.1874212598425
$co.Height = 324.0000787401575
